$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("HP")
$ws2 = $wb.Worksheets.Item("APPLE")

# Make sure "HP" is the active / selected sheet (it already is tabSelected in the source file).
$ws1.Activate()

# ---------------------------------------------------------------------------
# 1) Body rows (2-58): apply the thin-border-all-sides style (no bold, no fill)
#    to the full A:D block, including columns B and D that previously had no
#    cell entries at all. We reuse the already-existing border-only style
#    from the "APPLE" sheet (cell A2) via a format-only paste so the existing
#    style index is reused instead of minting a fresh one.
# ---------------------------------------------------------------------------
$ws2.Range("A2").Copy()
$ws1.Range("A2:D58").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) New row 59: fully empty cells A:D but still carrying the bordered style,
#    matching the pattern used for the rest of the table body.
# ---------------------------------------------------------------------------
$ws2.Range("A2").Copy()
$ws1.Range("A59:D59").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3) Header row 1, columns A-C: bordered + bold (no center alignment). Start
#    from the border-only style, then layer Bold on top so we don't disturb
#    alignment/number-format and keep reusing existing style plumbing as much
#    as possible.
# ---------------------------------------------------------------------------
$ws2.Range("A2").Copy()
$ws1.Range("A1:C1").PasteSpecial(-4122)
$ws1.Range("A1:C1").Font.Bold = $true

# ---------------------------------------------------------------------------
# 4) Header row 1, column D: hyperlink style + border (matches the existing
#    hyperlink-with-border style already used on the "APPLE" sheet's D1).
# ---------------------------------------------------------------------------
$ws2.Range("D1").Copy()
$ws1.Range("D1").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 5) View state: zoom + selected cell.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 95
[void]$ws1.Range("F14").Select()
